$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto prices / volume percentages scraped on Thu Dec  7 11:44:57 UTC 2023
$ws.Range("D2").Value = "43.302.99"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "2.237.63"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'230.06"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'0.639"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("D7").Value = "'63.76"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").Value = "'0.0947"
$ws.Range("E10").Value = "  -5.74%  "
$ws.Range("D11").Value = "'56.40"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "'26.60"
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").Value = "2.571.95"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'15.13"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "'5.99"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "'0.820"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "2.238.15"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "43.185.39"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").Value = "0.0₃0959"
$ws.Range("E20").Value = "  -5.34%  "
$ws.Range("D21").Value = "'72.83"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "'6.03"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "'245.47"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +32.90%  "
$ws.Range("D26").Value = "'2.40"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "'174.15"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").Value = "'9.66"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").Value = "'21.55"
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("E31").Value = "  -5.42%  "
$ws.Range("D32").Value = "'1.40"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "'4.88"
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("D35").Value = "'0.0672"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  -7.56%  "
$ws.Range("D38").Value = "'6.29"
$ws.Range("E38").Value = "  -5.43%  "
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").Value = "'0.0248"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").Value = "'8.57"
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("D44").Value = "'16.90"
$ws.Range("E44").Value = "  -2.84%  "
$ws.Range("D45").Value = "'96.04"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'0.000206"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").Value = "1.424.64"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "'9.81"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.74"
$ws.Range("E51").Value = "  +0.56%  "
